$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 5: StudyFilesTab script, mirroring the pattern of rows 2-4
$ws.Range("A5").Value = "StudyFilesTab"

$ws.Range("B5").Value = " MATCH (f:file)-->(s:study)`nMATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['COTC022'] and demo.breed in ['American Staffordshire Terrier','Mixed Breed'] and diag.primary_disease_site in ['Bone (Appendicular)']`nWITH DISTINCT f, s`nRETURN `n  coalesce(f.file_name, '') AS ``File Name``,`n  coalesce(f.file_type, '') AS ``File Type``,`n  coalesce(""study"", '') AS ``Association``,`n  coalesce(f.file_description, '') AS ``Description``,`n  coalesce(f.file_format, '') AS ``File Format``,`n  coalesce(f.file_size, '') AS ``Size``,`n  coalesce(s.clinical_study_designation,'') AS ``Study Code``"
$ws.Range("B5").WrapText = $true

$ws.Range("C5").Value = "MATCH (s:study)`n  MATCH (demo:demographic) `n  MATCH (diag:diagnosis)`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n`tWHERE s.clinical_study_designation IN ['COTC022'] and demo.breed in ['American Staffordshire Terrier','Mixed Breed'] and diag.primary_disease_site in ['Bone (Appendicular)']`n    `nOPTIONAL MATCH (s)<-[:member_of]-(c:case)`nOPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files , `n`tcount(DISTINCT(samp)) as number_of_sample , `n`tcount(DISTINCT(c.case_id)) as number_of_cases , `n`tcount(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("C5").WrapText = $true

$ws.Range("D5").Value = "TC01_Canine_StudyCOTC022-Breed_Diagnosis_PrimDiseaseSite_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC01_Canine_StudyCOTC022-Breed_Diagnosis_PrimDiseaseSite_WebData.xlsx"

# Match the updated view state: zoom to 70% and move the active selection to B2
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("B2").Select()
